$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (quote-prefix via NumberFormat) for D-column values that
# would otherwise be auto-parsed as numbers by Excel, so they keep their exact
# original text formatting (e.g. trailing zeros like "0.480").
$textCells = @("D5", "D6", "D10", "D12", "D14", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D36", "D37", "D39", "D40", "D41", "D43", "D46", "D47", "D48", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated coin values scraped by the GitHub Actions job.
$ws.Range("D2").Value = "67.030.94"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "3.120.45"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "577.18"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "174.07"
$ws.Range("E6").Value = "  +3.69%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.115.52"
$ws.Range("E8").Value = "  +1.33%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "6.43"
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "0.480"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "37.30"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "3.633.05"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").Value = "67.028.04"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("D18").Value = "7.13"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").Value = "3.118.10"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").Value = "16.27"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "478.11"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").Value = "0.713"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "7.69"
$ws.Range("E23").Value = "  +2.96%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "13.34"
$ws.Range("E24").Value = "  +3.22%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "83.97"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "10.01"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("D29").Value = "7.99"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").Value = "2.42"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("D33").Value = "0.0₃0975"
$ws.Range("E33").Value = "  -4.18%  "
$ws.Range("E34").Value = "  -2.32%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").Value = "5.87"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "0.977"
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").Value = "2.09"
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("D40").Value = "50.05"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").Value = "0.311"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").Value = "8.62"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").Value = "2.811.66"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("E45").Value = "  -10.23%  "
$ws.Range("D46").Value = "0.0356"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").Value = "380.10"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").Value = "135.91"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "24.84"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("E51").Value = "  -0.50%  "
